$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new data to row 20: D20 = 0.5, E20 = note text (with red font style, like other note cells)
$ws.Range("D20").Value = 0.5
$ws.Range("E20").Value = "added one to the draft class"
$ws.Range("E20").Font.Color = $ws.Range("E21").Font.Color

# Update selection on the sheet to match new active cell
$ws.Range("E33").Select()

